$wb = $excel.ActiveWorkbook

# A new "Not in need" dimension column is inserted as the first data column
# (C) on every sheet. The old "not falling within the PiN dimensions" data
# (which used to be the last data column, F) is in fact the data for this
# new "Not in need" column, so it moves into C, and the columns that used to
# sit in C/D/E (access / aggravating circumstances / learning condition)
# shift right by one into D/E/F. G (protected environment) is unaffected.

$ws = $wb.Worksheets.Item("idp")
$ws.Range("C1").Value = "Not in need"
$ws.Range("D1").Value = "access"
$ws.Range("E1").Value = "aggravating circumstances"
$ws.Range("F1").Value = "learning condition"
$ws.Range("C2").Value = 0.7106540107189164
$ws.Range("D2").Value = 0.190710457186663
$ws.Range("E2").Value = 0.09531458339445195
$ws.Range("F2").Value = 0.00332094869996863
$ws.Range("G2").Value = 0
$ws.Range("C3").Value = 0.3290499252879424
$ws.Range("D3").Value = 0.4673340598861581
$ws.Range("E3").Value = 0.02351335102101421
$ws.Range("F3").Value = 0.1363377937635117
$ws.Range("G3").Value = 0.0437648700413738
$ws.Range("C4").Value = 0.451490579415013
$ws.Range("D4").Value = 0.3029933244222449
$ws.Range("E4").Value = 0.02115781288198968
$ws.Range("F4").Value = 0.1668594842488218
$ws.Range("G4").Value = 0.05749879903193069
$ws.Range("C5").Value = 0.2418234127364369
$ws.Range("D5").Value = 0.6520581443625028
$ws.Range("E5").Value = 0.01366262932702547
$ws.Range("F5").Value = 0.07739927371974965
$ws.Range("G5").Value = 0.01505653985428526
$ws.Range("C6").Value = 0.03115927001577109
$ws.Range("D6").Value = 0.8414055018184412
$ws.Range("E6").Value = 0.006066203886073819
$ws.Range("F6").Value = 0.01720586734079309
$ws.Range("G6").Value = 0.1041631569389209
$ws.Range("C7").Value = 0.3453075914522154
$ws.Range("D7").Value = 0.5532912991791453
$ws.Range("E7").Value = 0.02136518339723263
$ws.Range("F7").Value = 0.01254642344682841
$ws.Range("G7").Value = 0.06748950252457825
$ws.Range("C8").Value = 0.6672571779664611
$ws.Range("D8").Value = 0.2367263812424759
$ws.Range("E8").Value = 0.06125870664438033
$ws.Range("F8").Value = 0.02324762293162015
$ws.Range("G8").Value = 0.01151011121506259
$ws.Range("C9").Value = 0.4300183427740319
$ws.Range("D9").Value = 0.4386203619332983
$ws.Range("E9").Value = 0.06879676346048313
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0.06256453183218642
$ws.Range("C10").Value = 0.503759477991108
$ws.Range("D10").Value = 0.270105571906739
$ws.Range("E10").Value = 0.06386966279094733
$ws.Range("F10").Value = 0.05680002453431085
$ws.Range("G10").Value = 0.1054652627768949
$ws.Range("C11").Value = 0.3042856291533364
$ws.Range("D11").Value = 0.5100784972268301
$ws.Range("E11").Value = 0.005852934652471125
$ws.Range("F11").Value = 0.02119727645262988
$ws.Range("G11").Value = 0.1585856625147326
$ws.Range("C12").Value = 0.5351307850709578
$ws.Range("D12").Value = 0.4023060687369476
$ws.Range("E12").Value = 0.03702069958023028
$ws.Range("F12").Value = 0.0255424466118643
$ws.Range("G12").Value = 0
$ws.Range("C13").Value = 0.1606215520728562
$ws.Range("D13").Value = 0.780119179566651
$ws.Range("E13").Value = 0.02185734521386637
$ws.Range("F13").Value = 0.03031483260684526
$ws.Range("G13").Value = 0.007087090539781185

$ws = $wb.Worksheets.Item("ndsp")
$ws.Range("C1").Value = "Not in need"
$ws.Range("D1").Value = "access"
$ws.Range("E1").Value = "aggravating circumstances"
$ws.Range("F1").Value = "learning condition"
$ws.Range("C2").Value = 0.3202038642640018
$ws.Range("D2").Value = 0.5997279723345624
$ws.Range("E2").Value = 0.01887686599956911
$ws.Range("F2").Value = 0.009269706648321113
$ws.Range("G2").Value = 0.05192159075354553

$ws = $wb.Worksheets.Item("ocap")
$ws.Range("C1").Value = "Not in need"
$ws.Range("D1").Value = "access"
$ws.Range("E1").Value = "aggravating circumstances"
$ws.Range("F1").Value = "learning condition"
$ws.Range("C2").Value = 0.7355010475254173
$ws.Range("D2").Value = 0.2246645324560284
$ws.Range("E2").Value = 0.02317121176486182
$ws.Range("F2").Value = 0.01006463778453275
$ws.Range("G2").Value = 0.006598570469159679
$ws.Range("C3").Value = 0.2881114296960966
$ws.Range("D3").Value = 0.6370232145138697
$ws.Range("E3").Value = 0.003355230387105165
$ws.Range("F3").Value = 0.03677214735238916
$ws.Range("G3").Value = 0.0347379780505395
$ws.Range("C4").Value = 0.542864226316311
$ws.Range("D4").Value = 0.3836873839858633
$ws.Range("E4").Value = 0.04386467216042554
$ws.Range("F4").Value = 0.0139614044321047
$ws.Range("G4").Value = 0.01562231310529551
$ws.Range("C5").Value = 0.2957402131124322
$ws.Range("D5").Value = 0.6475361562386643
$ws.Range("E5").Value = 0.02662675445858609
$ws.Range("F5").Value = 0.01034372974580771
$ws.Range("G5").Value = 0.0197531464445096
$ws.Range("C6").Value = 0.2168995155835067
$ws.Range("D6").Value = 0.7067734103739008
$ws.Range("E6").Value = 0.0479594500977146
$ws.Range("F6").Value = 0.02380464262637988
$ws.Range("G6").Value = 0.004562981318497784
$ws.Range("C7").Value = 0.7255237755561456
$ws.Range("D7").Value = 0.2031253855922261
$ws.Range("E7").Value = 0.03740669864710201
$ws.Range("F7").Value = 0.0339441402045264
$ws.Range("G7").Value = 0
$ws.Range("C8").Value = 0.7749837497270198
$ws.Range("D8").Value = 0.2011021168930733
$ws.Range("E8").Value = 0.01010426870802466
$ws.Range("F8").Value = 0.008707118220853691
$ws.Range("G8").Value = 0.005102746451028845
$ws.Range("C9").Value = 0.7567306197548092
$ws.Range("D9").Value = 0.1048223149121534
$ws.Range("E9").Value = 0.07569739695083222
$ws.Range("F9").Value = 0.06274966838220511
$ws.Range("G9").Value = 0
$ws.Range("C10").Value = 0.7815722567773472
$ws.Range("D10").Value = 0.2048177817451928
$ws.Range("E10").Value = 0.0136099614774599
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("C11").Value = 0.7738342948998559
$ws.Range("D11").Value = 0.1046527264059715
$ws.Range("E11").Value = 0.03543814744568882
$ws.Range("F11").Value = 0.0788805680596739
$ws.Range("G11").Value = 0.007194263188809723
$ws.Range("C12").Value = 0.6759834838293406
$ws.Range("D12").Value = 0.1111364861551428
$ws.Range("E12").Value = 0.07403526563460644
$ws.Range("F12").Value = 0.04505024459947984
$ws.Range("G12").Value = 0.0937945197814303
$ws.Range("C13").Value = 0.4619138320502865
$ws.Range("D13").Value = 0.403750898051613
$ws.Range("E13").Value = 0.01966706479867863
$ws.Range("F13").Value = 0.07811591978181062
$ws.Range("G13").Value = 0.03655228531761142
$ws.Range("C14").Value = 0.6577904324554799
$ws.Range("D14").Value = 0.2580592382461159
$ws.Range("E14").Value = 0.07780282150174837
$ws.Range("F14").Value = 0.002059391714946343
$ws.Range("G14").Value = 0.004288116081709648
$ws.Range("C15").Value = 0.8174822877535954
$ws.Range("D15").Value = 0.07499063577569208
$ws.Range("E15").Value = 0.04406802745667497
$ws.Range("F15").Value = 0.06345904901403755
$ws.Range("G15").Value = 0
$ws.Range("C16").Value = 0.538121125714962
$ws.Range("D16").Value = 0.3438485809116489
$ws.Range("E16").Value = 0.05173825925298021
$ws.Range("F16").Value = 0.05389118409563317
$ws.Range("G16").Value = 0.01240085002477574
$ws.Range("C17").Value = 0.8548534137558078
$ws.Range("D17").Value = 0.1094267861545242
$ws.Range("E17").Value = 0.03571980008966798
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("C18").Value = 0.9248139239710328
$ws.Range("D18").Value = 0.03634751776327397
$ws.Range("E18").Value = 0.03292999288485467
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0.005908565380838581
$ws.Range("C19").Value = 0.87123704802994
$ws.Range("D19").Value = 0.08039684079360003
$ws.Range("E19").Value = 0.04347052377995093
$ws.Range("F19").Value = 0.004895587396509143
$ws.Range("G19").Value = 0

$ws = $wb.Worksheets.Item("ret")
$ws.Range("C1").Value = "Not in need"
$ws.Range("D1").Value = "access"
$ws.Range("E1").Value = "aggravating circumstances"
$ws.Range("F1").Value = "learning condition"
$ws.Range("C2").Value = 0.5993405930558748
$ws.Range("D2").Value = 0.3461159896987495
$ws.Range("E2").Value = 0.03700134904572487
$ws.Range("F2").Value = 0.01754206819965072
$ws.Range("G2").Value = 0
$ws.Range("C3").Value = 0.1643234671516556
$ws.Range("D3").Value = 0.5363990822578679
$ws.Range("E3").Value = 0.01688281685020779
$ws.Range("F3").Value = 0.1905420565106229
$ws.Range("G3").Value = 0.09185257722964593
$ws.Range("C4").Value = 0.4917716660792544
$ws.Range("D4").Value = 0.4861326132088599
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0.01567743048028169
$ws.Range("G4").Value = 0.006418290231603985
$ws.Range("C5").Value = 0.4068563174242846
$ws.Range("D5").Value = 0.5296564849408226
$ws.Range("E5").Value = 0.03816465636136621
$ws.Range("F5").Value = 0.01414738013903636
$ws.Range("G5").Value = 0.01117516113449035
$ws.Range("C6").Value = 0.07562392814612101
$ws.Range("D6").Value = 0.9193078839006363
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0.005068187953242716
$ws.Range("G6").Value = 0
$ws.Range("C7").Value = 0.5497172552716943
$ws.Range("D7").Value = 0.3892670239583827
$ws.Range("E7").Value = 0.06101572076992296
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("C8").Value = 0.4732055462978784
$ws.Range("D8").Value = 0.4160086106878902
$ws.Range("E8").Value = 0.0432878883621454
$ws.Range("F8").Value = 0.06749795465208604
$ws.Range("G8").Value = 0
$ws.Range("C9").Value = 0.2802404301971403
$ws.Range("D9").Value = 0.4304411167262168
$ws.Range("E9").Value = 0.08493390733218803
$ws.Range("F9").Value = 0.1063556480784028
$ws.Range("G9").Value = 0.09802889766605188
$ws.Range("C10").Value = 0.5337870602633682
$ws.Range("D10").Value = 0.2525549026166085
$ws.Range("E10").Value = 0.1393912851416745
$ws.Range("F10").Value = 0.05695902476239229
$ws.Range("G10").Value = 0.0173077272159565
$ws.Range("C11").Value = 0.3287750150324628
$ws.Range("D11").Value = 0.5242081101445487
$ws.Range("E11").Value = 0.04803177498858425
$ws.Range("F11").Value = 0.05975679914582763
$ws.Range("G11").Value = 0.03922830068857657
$ws.Range("C12").Value = 0.5928738109452726
$ws.Range("D12").Value = 0.1980015235867318
$ws.Range("E12").Value = 0.2006406549106493
$ws.Range("F12").Value = 0.008484010557346435
$ws.Range("G12").Value = 0
$ws.Range("C13").Value = 0.3666693062507451
$ws.Range("D13").Value = 0.4897918112166131
$ws.Range("E13").Value = 0.06306542781760054
$ws.Range("F13").Value = 0.02701275396855353
$ws.Range("G13").Value = 0.0534607007464878
